# Insert a new data row before existing row 6, shifting all subsequent
# rows (previously 6..168) down to 7..169, and populate the new row 6
# with the latest week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this shifts row 6..168 down to 7..169
# and carries formatting (e.g. the date style on column D) along with it.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = Get-Date -Year 2022 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 100112037
$ws.Cells.Item(6, 7).Value = "Cebollín"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 600
$ws.Cells.Item(6, 11).Value = 1100
$ws.Cells.Item(6, 12).Value = 1200
$ws.Cells.Item(6, 13).Value = 1150
$ws.Cells.Item(6, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(6, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(6, 16).Value = 192
$ws.Cells.Item(6, 17).Value = 6
$ws.Cells.Item(6, 18).Value = "Hortaliza"
